$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title: merge the two runs ("Hydr" + bookmark + "oxyproline Assay")
#    into a single run "Hydroxyproline Assay". The _GoBack bookmark
#    that used to sit between them is removed here (it gets re-added
#    later, at its new location near the end of the "grind the
#    tissue." paragraph).
# ------------------------------------------------------------------
$d.Content.Find.Execute("Hydroxyproline Assay", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Hydroxyproline Assay", 2)

# ------------------------------------------------------------------
# 2. Re-insert the _GoBack bookmark right after "...to properly grind
#    the tissue." (collapsed, i.e. bookmarkStart immediately followed
#    by bookmarkEnd, with no text in between). A truly zero-length
#    Range confuses this host's Bookmarks.Add, so we briefly insert a
#    one-character placeholder, wrap the bookmark around it, then
#    delete the placeholder -- leaving the bookmark collapsed exactly
#    where we need it.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("to properly grind the tissue.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$endPos = $anchor.End

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$markerRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$cleanup = $d.Range($endPos, $endPos + 1)
$cleanup.Text = ""

# ------------------------------------------------------------------
# 3. Move the word "turn " from the plain run to the start of the
#    highlighted/glow run: ". You will want to turn " + "this on now"
#    -> ". You will want to " + "turn this on now"
# ------------------------------------------------------------------
$d.Content.Find.Execute(". You will want to turn ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ". You will want to ", 2)
$d.Content.Find.Execute("this on now", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "turn this on now", 2)
